$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three vendors that are being dropped entirely from the map.
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows(79).Delete()   # Wood Stove Kitchen
$ws.Rows(54).Delete()   # Painting as Art & Ritual
$ws.Rows(44).Delete()   # Mugger's Marrow, LLC

# Renumber the Index column (A) sequentially now that rows were removed.
for ($r = 2; $r -le 76; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Fill in the Zip Code column (C) for every vendor row, in the same
# order the cells were originally populated (rows 75 and 74 were
# entered out of sequence relative to the rest of the column).
$rowOrder = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,75,74,76)
$zipByRow = @{
    2="05647"; 3="02907"; 4="03304"; 5="02116"; 6="02467"; 7="01852"; 8="02139"; 9="01915"; 10="02675";
    11="02842"; 12="01230"; 13="02126"; 14="02452"; 15="02139"; 16="02129"; 17="01801"; 18="02120"; 19="02143";
    20="02576"; 21="02653"; 22="02121"; 23="01027"; 24="02111"; 25="02136"; 26="02149"; 27="01776"; 28="01915";
    29="04856"; 30="02906"; 31="02481"; 32="01375"; 33="01752"; 34="02109"; 35="03820"; 36="02823"; 37="02151";
    38="01035"; 39="02114"; 40="02461"; 41="02118"; 42="01984"; 43="02131"; 44="02114"; 45="01742"; 46="02601";
    47="02138"; 48="02115"; 49="02113"; 50="01002"; 51="01702"; 52="01301"; 53="02457"; 54="02144"; 55="02131";
    56="01301"; 57="01331"; 58="01966"; 59="05482"; 60="02143"; 61="05701"; 62="01915"; 63="02144"; 64="02885";
    65="02143"; 66="02143"; 67="02081"; 68="02356"; 69="02633"; 70="02110"; 71="01801"; 72="01983"; 73="03263";
    74="10956"; 75="02043"; 76="20170"
}

foreach ($r in $rowOrder) {
    $ws.Cells.Item($r, 3).Value = $zipByRow[$r]
}

# Update the selected cell, mirroring the author's last selection in Excel.
$ws.Range("B21").Select()
